$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.301.65"
$ws.Range("E2").Value = "  -1.02%  "

$ws.Range("D3").Value = "3.534.16"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "608.13"
$ws.Range("E5").Value = "  +0.53%  "

$ws.Range("D6").Value = "145.15"
$ws.Range("E6").Value = "  -2.02%  "

$ws.Range("D7").Value = "3.536.64"
$ws.Range("E7").Value = "  +0.62%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  -0.38%  "

$ws.Range("D10").Value = "0.137"
$ws.Range("E10").Value = "  -4.23%  "

$ws.Range("E11").Value = "  +2.15%  "

$ws.Range("E12").Value = "  -2.16%  "

$ws.Range("D13").Value = "4.130.45"
$ws.Range("E13").Value = "  +0.45%  "

$ws.Range("D14").Value = "0.0000209"
$ws.Range("E14").Value = "  -3.09%  "

$ws.Range("D15").Value = "30.50"
$ws.Range("E15").Value = "  -3.55%  "

$ws.Range("D16").Value = "3.526.72"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").Value = "66.347.27"
$ws.Range("E17").Value = "  -0.99%  "

$ws.Range("E18").Value = "  +0.05%  "

$ws.Range("D19").Value = "10.89"
$ws.Range("E19").Value = "  +1.75%  "

$ws.Range("D20").Value = "6.22"
$ws.Range("E20").Value = "  -2.70%  "

$ws.Range("D21").Value = "14.96"
$ws.Range("E21").Value = "  -2.67%  "

$ws.Range("D22").Value = "427.24"
$ws.Range("E22").Value = "  -1.91%  "

$ws.Range("D23").Value = "0.602"
$ws.Range("E23").Value = "  -1.37%  "

$ws.Range("D24").Value = "78.49"
$ws.Range("E24").Value = "  -1.33%  "

$ws.Range("D25").Value = "3.664.85"
$ws.Range("E25").Value = "  +0.28%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").Value = "0.0000120"
$ws.Range("E27").Value = "  -0.58%  "

$ws.Range("D28").Value = "9.32"
$ws.Range("E28").Value = "  -5.36%  "

$ws.Range("D29").Value = "8.04"
$ws.Range("E29").Value = "  -3.15%  "

$ws.Range("E30").Value = "  -1.13%  "

$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("E32").Value = "  -1.88%  "

$ws.Range("E33").Value = "  -7.16%  "

$ws.Range("D34").Value = "25.34"
$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("D35").Value = "3.517.36"
$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  -3.44%  "

$ws.Range("E38").Value = "  -4.69%  "

$ws.Range("E39").Value = "  -2.67%  "

$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.23%  "

$ws.Range("D41").Value = "170.72"
$ws.Range("E41").Value = "  +0.80%  "

$ws.Range("D42").Value = "0.0861"
$ws.Range("E42").Value = "  -3.49%  "

$ws.Range("E43").Value = "  -4.41%  "

$ws.Range("D44").Value = "0.892"
$ws.Range("E44").Value = "  -0.47%  "

$ws.Range("E45").Value = "  -9.95%  "

$ws.Range("D46").Value = "45.51"
$ws.Range("E46").Value = "  -0.43%  "

$ws.Range("D47").Value = "1.22"
$ws.Range("E47").Value = "  -7.41%  "

$ws.Range("D48").Value = "25.95"
$ws.Range("E48").Value = "  -10.26%  "

$ws.Range("D49").Value = "2.42"
$ws.Range("E49").Value = "  -0.88%  "

$ws.Range("D50").Value = "7.20"
$ws.Range("E50").Value = "  -3.73%  "

$ws.Range("D51").Value = "0.950"
$ws.Range("E51").Value = "  -3.77%  "
